$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update I2:I10 from 1 to 2, and clear J2:J10 (validation columns)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = 2
    $ws.Cells.Item($r, 10).Value = $null
}

# Sheet view: reset scroll position and move the selection to B14
$ws.Range("B14").Select()
